$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 151 ("coming soon" post), shifting subsequent rows up.
$ws.Rows.Item(151).Delete()
